$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header columns (row 1) ---
$ws.Cells.Item(1, 1).Value = "mx_state"
$ws.Cells.Item(1, 2).Value = "mx_municipality"
$ws.Cells.Item(1, 3).Value = "n_matriculas"
$ws.Cells.Item(1, 4).Value = "pct_matriculas"

# --- Title-case the connector words (de/del/el/la/los/las/y) in state & municipality names ---
$ws.Cells.Item(3, 2).Value = "Pabellón De Arteaga"
$ws.Cells.Item(13, 2).Value = "Amatenango De La Frontera"
$ws.Cells.Item(28, 2).Value = "Mazapa De Madero"
$ws.Cells.Item(53, 1).Value = "Ciudad De México"
$ws.Cells.Item(57, 2).Value = "Cuajimalpa De Morelos"
$ws.Cells.Item(74, 1).Value = "Estado De México"
$ws.Cells.Item(74, 2).Value = "Acambay De Ruíz Castañeda"
$ws.Cells.Item(81, 2).Value = "Atizapán De Zaragoza"
$ws.Cells.Item(86, 2).Value = "Chapa De Mota"
$ws.Cells.Item(89, 2).Value = "Coacalco De Berriozábal"
$ws.Cells.Item(93, 2).Value = "Ecatepec De Morelos"
$ws.Cells.Item(100, 2).Value = "Naucalpan De Juárez"
$ws.Cells.Item(105, 2).Value = "San Felipe Del Progreso"
$ws.Cells.Item(115, 2).Value = "Tlalnepantla De Baz"
$ws.Cells.Item(123, 2).Value = "San Miguel De Allende"
$ws.Cells.Item(124, 2).Value = "Apaseo El Alto"
$ws.Cells.Item(125, 2).Value = "Apaseo El Grande"
$ws.Cells.Item(129, 2).Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Cells.Item(131, 2).Value = "Jaral Del Progreso"
$ws.Cells.Item(134, 2).Value = "Purísima Del Rincón"
$ws.Cells.Item(137, 2).Value = "San Diego De La Unión"
$ws.Cells.Item(145, 2).Value = "Acapulco De Juárez"
$ws.Cells.Item(148, 2).Value = "Alcozauca De Guerrero"
$ws.Cells.Item(152, 2).Value = "Atoyac De Álvarez"
$ws.Cells.Item(153, 2).Value = "Ayutla De Los Libres"
$ws.Cells.Item(154, 2).Value = "Buenavista De Cuéllar"
$ws.Cells.Item(155, 2).Value = "Chilapa De Álvarez"
$ws.Cells.Item(156, 2).Value = "Chilpancingo De Los Bravo"
$ws.Cells.Item(157, 2).Value = "Coahuayutla De José María Izazaga"
$ws.Cells.Item(160, 2).Value = "Coyuca De Benítez"
$ws.Cells.Item(161, 2).Value = "Coyuca De Catalán"
$ws.Cells.Item(163, 2).Value = "Cutzamala De Pinzón"
$ws.Cells.Item(167, 2).Value = "Iguala De La Independencia"
$ws.Cells.Item(169, 2).Value = "Ixcateopan De Cuauhtémoc"
$ws.Cells.Item(170, 2).Value = "Zihuatanejo De Azueta"
$ws.Cells.Item(172, 2).Value = "La Unión De Isidoro Montes De Oca"
$ws.Cells.Item(174, 2).Value = "Mártir De Cuilapan"
$ws.Cells.Item(183, 2).Value = "Taxco De Alarcón"
$ws.Cells.Item(186, 2).Value = "Tepecoacuilco De Trujano"
$ws.Cells.Item(188, 2).Value = "Tixtla De Guerrero"
$ws.Cells.Item(190, 2).Value = "Tlalixtaquilla De Maldonado"
$ws.Cells.Item(191, 2).Value = "Tlapa De Comonfort"
$ws.Cells.Item(202, 2).Value = "Atotonilco El Grande"
$ws.Cells.Item(206, 2).Value = "Huejutla De Reyes"
$ws.Cells.Item(211, 2).Value = "Mixquiahuala De Juárez"
$ws.Cells.Item(212, 2).Value = "Molango De Escamilla"
$ws.Cells.Item(214, 2).Value = "Omitlán De Juárez"
$ws.Cells.Item(215, 2).Value = "Pachuca De Soto"
$ws.Cells.Item(217, 2).Value = "Progreso De Obregón"
$ws.Cells.Item(219, 2).Value = "Santiago De Anaya"
$ws.Cells.Item(220, 2).Value = "Tenango De Doria"
$ws.Cells.Item(222, 2).Value = "Tepeji Del Río De Ocampo"
$ws.Cells.Item(226, 2).Value = "Tula De Allende"
$ws.Cells.Item(227, 2).Value = "Tulancingo De Bravo"
$ws.Cells.Item(231, 2).Value = "Ahualulco De Mercado"
$ws.Cells.Item(239, 2).Value = "Huejuquilla El Alto"
$ws.Cells.Item(242, 2).Value = "Jilotlán De Los Dolores"
$ws.Cells.Item(244, 2).Value = "Lagos De Moreno"
$ws.Cells.Item(250, 2).Value = "San Juan De Los Lagos"
$ws.Cells.Item(252, 2).Value = "Tamazula De Gordiano"
$ws.Cells.Item(253, 2).Value = "Tepatitlán De Morelos"
$ws.Cells.Item(255, 2).Value = "Tizapán El Alto"
$ws.Cells.Item(257, 2).Value = "Valle De Juárez"
$ws.Cells.Item(260, 2).Value = "Zapotlán El Grande"
$ws.Cells.Item(301, 2).Value = "Jonacatepec De Leandro Valle"
$ws.Cells.Item(304, 2).Value = "Puente De Ixtla"
$ws.Cells.Item(309, 2).Value = "Tetela Del Volcán"
$ws.Cells.Item(316, 2).Value = "Zacualpan De Amilpas"
$ws.Cells.Item(326, 2).Value = "Acatlán De Pérez Figueroa"
$ws.Cells.Item(328, 2).Value = "Ayoquezco De Aldama"
$ws.Cells.Item(329, 2).Value = "Ciénega De Zimatlán"
$ws.Cells.Item(332, 2).Value = "Fresnillo De Trujano"
$ws.Cells.Item(333, 2).Value = "Heroica Ciudad De Huajuapan De León"
$ws.Cells.Item(334, 2).Value = "Heroica Ciudad De Tlaxiaco"
$ws.Cells.Item(335, 2).Value = "Ixtlán De Juárez"
$ws.Cells.Item(336, 2).Value = "Mariscala De Juárez"
$ws.Cells.Item(337, 2).Value = "Oaxaca De Juárez"
$ws.Cells.Item(338, 2).Value = "Ocotlán De Morelos"
$ws.Cells.Item(339, 2).Value = "Putla Villa De Guerrero"
$ws.Cells.Item(345, 2).Value = "San Antonino El Alto"
$ws.Cells.Item(363, 2).Value = "San Miguel Del Puerto"
$ws.Cells.Item(368, 2).Value = "San Pedro Y San Pablo Teposcolula"
$ws.Cells.Item(374, 2).Value = "Santa Inés Del Monte"
$ws.Cells.Item(391, 2).Value = "Teotitlán De Flores Magón"
$ws.Cells.Item(392, 2).Value = "Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca"
$ws.Cells.Item(393, 2).Value = "Tlacolula De Matamoros"
$ws.Cells.Item(394, 2).Value = "Tlalixtac De Cabrera"
$ws.Cells.Item(395, 2).Value = "Totontepec Villa De Morelos"
$ws.Cells.Item(398, 2).Value = "Villa De Tututepec"
$ws.Cells.Item(399, 2).Value = "Villa De Zaachila"
$ws.Cells.Item(402, 2).Value = "Zimatlán De Álvarez"
$ws.Cells.Item(424, 2).Value = "Chalchicomula De Sesma"
$ws.Cells.Item(429, 2).Value = "Chila De La Sal"
$ws.Cells.Item(437, 2).Value = "Cuayuca De Andrade"
$ws.Cells.Item(438, 2).Value = "Cuetzalan Del Progreso"
$ws.Cells.Item(447, 2).Value = "Huehuetlán El Chico"
$ws.Cells.Item(448, 2).Value = "Huehuetlán El Grande"
$ws.Cells.Item(452, 2).Value = "Izúcar De Matamoros"
$ws.Cells.Item(459, 2).Value = "Los Reyes De Juárez"
$ws.Cells.Item(466, 2).Value = "Palmar De Bravo"
$ws.Cells.Item(485, 2).Value = "San Nicolás De Los Ranchos"
$ws.Cells.Item(488, 2).Value = "San Salvador El Seco"
$ws.Cells.Item(489, 2).Value = "San Salvador El Verde"
$ws.Cells.Item(495, 2).Value = "Tecali De Herrera"
$ws.Cells.Item(502, 2).Value = "Tepanco De López"
$ws.Cells.Item(503, 2).Value = "Tepatlaxco De Hidalgo"
$ws.Cells.Item(508, 2).Value = "Tepexi De Rodríguez"
$ws.Cells.Item(510, 2).Value = "Tetela De Ocampo"
$ws.Cells.Item(515, 2).Value = "Tlacotepec De Benito Juárez"
$ws.Cells.Item(527, 2).Value = "Xayacatlán De Bravo"
$ws.Cells.Item(545, 2).Value = "San Juan Del Río"
$ws.Cells.Item(552, 2).Value = "Ciudad Del Maíz"
$ws.Cells.Item(554, 2).Value = "Mexquitic De Carmona"
$ws.Cells.Item(558, 2).Value = "Santa María Del Río"
$ws.Cells.Item(560, 2).Value = "Tanquián De Escobedo"
$ws.Cells.Item(576, 2).Value = "Jalpa De Méndez"
$ws.Cells.Item(589, 2).Value = "Amaxac De Guerrero"
$ws.Cells.Item(596, 2).Value = "Ixtacuixtla De Mariano Matamoros"
$ws.Cells.Item(598, 2).Value = "Mazatecochco De José María Morelos"
$ws.Cells.Item(601, 2).Value = "Papalotla De Xicohténcatl"
$ws.Cells.Item(605, 2).Value = "San Pablo Del Monte"
$ws.Cells.Item(606, 2).Value = "Sanctórum De Lázaro Cárdenas"
$ws.Cells.Item(611, 2).Value = "Tepetitla De Lardizábal"
$ws.Cells.Item(626, 2).Value = "Alto Lucero De Gutiérrez Barrios"
$ws.Cells.Item(629, 2).Value = "Boca Del Río"
$ws.Cells.Item(638, 2).Value = "Cosamaloapan De Carpio"
$ws.Cells.Item(642, 2).Value = "Hueyapan De Ocampo"
$ws.Cells.Item(643, 2).Value = "Ignacio De La Llave"
$ws.Cells.Item(644, 2).Value = "Ixhuatlán De Madero"
$ws.Cells.Item(649, 2).Value = "Juchique De Ferrer"
$ws.Cells.Item(652, 2).Value = "Las Vigas De Ramírez"
$ws.Cells.Item(653, 2).Value = "Lerdo De Tejada"
$ws.Cells.Item(655, 2).Value = "Martínez De La Torre"
$ws.Cells.Item(661, 2).Value = "Paso De Ovejas"
$ws.Cells.Item(664, 2).Value = "Poza Rica De Hidalgo"
$ws.Cells.Item(669, 2).Value = "Soledad De Doblado"
$ws.Cells.Item(688, 2).Value = "Nochistlán De Mejía"
$ws.Cells.Item(690, 2).Value = "Villa De Cos"

# --- Fix two floating point rounding values ---
$ws.Cells.Item(60, 4).Value = 0.009218289085545724
$ws.Cells.Item(272, 4).Value = 0.009218289085545724

# --- Remove trailing footer/metadata rows 695-699 ---
$ws.Range("A695:D699").EntireRow.Delete()
